$d = $word.ActiveDocument

# Locate the paragraph whose entire text is "First" (the answer to
# "How many classes have you completed?"). The new Q&A pair
# ("Where are you located?" / "California") must be inserted right after it.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -eq "First`r") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ge 1) {

    # --- New bold question paragraph: "Where are you located?" ---
    # Copy an existing bold Q-style paragraph (paragraph 1, "Which Program
    # are you enrolled in?") so the new paragraph picks up the exact same
    # bold run/paragraph-mark formatting (b + bCs), then paste it right
    # after the "First" paragraph and overwrite its text.
    $boldSourcePara = $d.Paragraphs(1)
    $boldSourceRange = $d.Range($boldSourcePara.Range.Start, $boldSourcePara.Range.End)
    $boldSourceRange.Copy()

    $firstPara = $d.Paragraphs($targetIndex)
    $pasteAt = $firstPara.Range.End
    $pasteRange = $d.Range($pasteAt, $pasteAt)
    $pasteRange.Paste()

    $questionPara = $d.Paragraphs($targetIndex + 1)
    $questionTextRange = $d.Range($questionPara.Range.Start, $questionPara.Range.End - 1)
    $questionTextRange.Text = "Where are you located?"

    # --- New plain answer paragraph: "California" ---
    # Copy an existing plain (non-bold, no direct formatting) answer
    # paragraph (paragraph 2, "OMSCS") so the new paragraph has no rPr/pPr
    # at all, then paste it after the question paragraph and overwrite its
    # text.
    $plainSourcePara = $d.Paragraphs(2)
    $plainSourceRange = $d.Range($plainSourcePara.Range.Start, $plainSourcePara.Range.End)
    $plainSourceRange.Copy()

    $questionPara = $d.Paragraphs($targetIndex + 1)
    $pasteAt2 = $questionPara.Range.End
    $pasteRange2 = $d.Range($pasteAt2, $pasteAt2)
    $pasteRange2.Paste()

    $answerPara = $d.Paragraphs($targetIndex + 2)
    $answerTextRange = $d.Range($answerPara.Range.Start, $answerPara.Range.End - 1)
    $answerTextRange.Text = "California"
}
